$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full data block for rows 2-17 (header row 1 is unchanged):
# Columns A-D are the categorical labels (Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster); a 4th sending/target cluster ("ECs")
# has been added, expanding the table from 12 to 16 data rows.
$data = @(
    ,("ECs", "Ntn3", "Neo1", "ECs", 2, 0.6666666666666666, 0.3056536666666667, 0.916961, 0.09028452848253411, 0.09028452848253413, 3, 1, 2.373057666666667, 7.119173, 0.03663411841135517, 0.03663411841135517, 0.7253337770281111, 6.528003993253, 0.003307494107142523, 0.003307494107142523)
    ,("ECs", "Ntn3", "Neo1", "FAPs", 2, 0.6666666666666666, 0.3056536666666667, 0.916961, 0.09028452848253411, 0.09028452848253413, 3, 1, 43.26393666666667, 129.79181, 0.6678877639810287, 0.6678877639810287, 13.22378087660111, 119.01402788941, 0.0602999318502812, 0.06029993185028121)
    ,("ECs", "Ntn3", "Neo1", "M2", 2, 0.6666666666666666, 0.3056536666666667, 0.916961, 0.09028452848253411, 0.09028452848253413, 3, 1, 0.4664283333333333, 1.399285, 0.007200495391983466, 0.007200495391983466, 0.1425655303205555, 1.283089772885, 0.0006500933313058868, 0.0006500933313058869)
    ,("ECs", "Ntn3", "Neo1", "sCs", 2, 0.6666666666666666, 0.3056536666666667, 0.916961, 0.09028452848253411, 0.09028452848253413, 3, 1, 18.67383333333333, 56.0215, 0.2882776222156329, 0.2882776222156329, 5.707725629055556, 51.36953066150001, 0.02602700919380451, 0.02602700919380452)
    ,("FAPs", "Ntn3", "Neo1", "ECs", 3, 1, 2.160221, 6.480663, 0.6380899549808606, 0.6380899549808607, 3, 1, 2.373057666666667, 7.119173, 0.03663411841135517, 0.03663411841135517, 5.126329005744333, 46.136961051699, 0.02337586296786514, 0.02337586296786514)
    ,("FAPs", "Ntn3", "Neo1", "FAPs", 3, 1, 2.160221, 6.480663, 0.6380899549808606, 0.6380899549808607, 3, 1, 43.26393666666667, 129.79181, 0.6678877639810287, 0.6678877639810287, 93.45966453000332, 841.1369807700299, 0.4261724732509223, 0.4261724732509223)
    ,("FAPs", "Ntn3", "Neo1", "M2", 3, 1, 2.160221, 6.480663, 0.6380899549808606, 0.6380899549808607, 3, 1, 0.4664283333333333, 1.399285, 0.007200495391983466, 0.007200495391983466, 1.007588280661667, 9.068294525954999, 0.004594563780510624, 0.004594563780510625)
    ,("FAPs", "Ntn3", "Neo1", "sCs", 3, 1, 2.160221, 6.480663, 0.6380899549808606, 0.6380899549808607, 3, 1, 18.67383333333333, 56.0215, 0.2882776222156329, 0.2882776222156329, 40.33960691716667, 363.0564622545, 0.1839470549815627, 0.1839470549815627)
    ,("M2", "Ntn3", "Neo1", "ECs", 3, 1, 0.3550333333333333, 1.0651, 0.1048703830225572, 0.1048703830225572, 3, 1, 2.373057666666667, 7.119173, 0.03663411841135517, 0.03663411841135517, 0.8425145735888888, 7.582631162299999, 0.003841834029492531, 0.003841834029492531)
    ,("M2", "Ntn3", "Neo1", "FAPs", 3, 1, 0.3550333333333333, 1.0651, 0.1048703830225572, 0.1048703830225572, 3, 1, 43.26393666666667, 129.79181, 0.6678877639810287, 0.6678877639810287, 15.36013964788889, 138.241256831, 0.07004164562476976, 0.07004164562476978)
    ,("M2", "Ntn3", "Neo1", "M2", 3, 1, 0.3550333333333333, 1.0651, 0.1048703830225572, 0.1048703830225572, 3, 1, 0.4664283333333333, 1.399285, 0.007200495391983466, 0.007200495391983466, 0.1655976059444444, 1.4903784535, 0.0007551187097094642, 0.0007551187097094643)
    ,("M2", "Ntn3", "Neo1", "sCs", 3, 1, 0.3550333333333333, 1.0651, 0.1048703830225572, 0.1048703830225572, 3, 1, 18.67383333333333, 56.0215, 0.2882776222156329, 0.2882776222156329, 6.629833294444444, 59.66849965, 0.03023178465858546, 0.03023178465858547)
    ,("sCs", "Ntn3", "Neo1", "ECs", 3, 1, 0.564541, 1.693623, 0.1667551335140479, 0.1667551335140479, 3, 1, 2.373057666666667, 7.119173, 0.03663411841135517, 0.03663411841135517, 1.339688348197666, 12.057195133779, 0.00610892730685497, 0.00610892730685497)
    ,("sCs", "Ntn3", "Neo1", "FAPs", 3, 1, 0.564541, 1.693623, 0.1667551335140479, 0.1667551335140479, 3, 1, 43.26393666666667, 129.79181, 0.6678877639810287, 0.6678877639810287, 24.42426606973666, 219.81839462763, 0.1113737132550553, 0.1113737132550553)
    ,("sCs", "Ntn3", "Neo1", "M2", 3, 1, 0.564541, 1.693623, 0.1667551335140479, 0.1667551335140479, 3, 1, 0.4664283333333333, 1.399285, 0.007200495391983466, 0.007200495391983466, 0.2633179177283333, 2.369861259554999, 0.001200719570457489, 0.001200719570457489)
    ,("sCs", "Ntn3", "Neo1", "sCs", 3, 1, 0.564541, 1.693623, 0.1667551335140479, 0.1667551335140479, 3, 1, 18.67383333333333, 56.0215, 0.2882776222156329, 0.2882776222156329, 10.54214454383333, 94.87930089450001, 0.04807177338168012, 0.04807177338168012)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
